# Add "Velocity in km/s" column (D) with formula converting ft/s to km/s,
# and a hand-entered rounded copy of those values in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "Velocity in km/s"

# Column D: formula B*0.0003048*1000 for rows 2..10.
# D2 is entered on its own (regular formula); D3:D10 are entered as a
# separate range-formula assignment so Excel groups them as one shared-
# formula block, matching the recorded end-state where D2 stays ungrouped
# and D3:D10 share formula si="0".
$ws.Range("D2").Formula = "=B2*0.0003048*1000"
$ws.Range("D3:D10").Formula = "=B3*0.0003048*1000"

# Column E: rounded literal values (not formulas), pasted as values
$ws.Range("E2").Value = 0.6892431
$ws.Range("E3").Value = 1.1359399999999999
$ws.Range("E4").Value = 0.034021000000000003
$ws.Range("E5").Value = 0.49453000000000003
$ws.Range("E6").Value = 0.23743
$ws.Range("E7").Value = 0.59843000000000002
$ws.Range("E8").Value = 0.3
$ws.Range("E9").Value = 0.17721000000000001
$ws.Range("E10").Value = 0.13408999999999999

# E4 uses a scientific-notation number format (numFmtId 11) in the target file
$ws.Range("E4").NumberFormat = "0.00E+00"

# Update selection to match the recorded end-state
$ws.Range("F18").Select()
